# Plantilla Lista de Tareas de la 6ta Iteracion
# - Actualiza estatus de un par de tareas a "Hecho" (lo que deja sin uso el
#   texto compartido "En proceso")
# - Registra 2 horas consumidas el dia correspondiente a la columna T (fila 10)
# - Re-crea las celdas combinadas del encabezado (fila 4) en el nuevo orden
# - Actualiza la celda activa seleccionada en la hoja "Casos de Uso"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Actualizar estatus de tareas
$ws.Range("F6").Value2 = "Hecho"
$ws.Range("F10").Value2 = "Hecho"

# Registrar horas consumidas (columna "Cons." del dia correspondiente) en la fila 10
$ws.Range("T10").Value2 = 2

# Reconstruir las celdas combinadas del encabezado en el nuevo orden
$mergedRangesOriginalOrder = @("AL4:AM4","H4:I4","K4:L4","N4:O4","Q4:R4","T4:U4","W4:X4","Z4:AA4","AC4:AD4","AF4:AG4","AI4:AJ4","AZ4:BA4","AO4:AP4","AR4:AS4","AU4:AV4","AX4:AY4")
foreach ($r in $mergedRangesOriginalOrder) {
    $ws.Range($r).UnMerge()
}
$mergedRangesNewOrder = @("AZ4:BA4","AO4:AP4","AR4:AS4","AU4:AV4","AX4:AY4","AL4:AM4","H4:I4","K4:L4","N4:O4","Q4:R4","T4:U4","W4:X4","Z4:AA4","AC4:AD4","AF4:AG4","AI4:AJ4")
foreach ($r in $mergedRangesNewOrder) {
    $ws.Range($r).Merge()
}

# Actualizar celda activa / seleccion visible
$ws.Activate() | Out-Null
$ws.Range("K9").Select() | Out-Null
